# Update FuelPrices at 2025-03-19 14:59
# - Row 5's Date cell (B5) switches from the "date-only" format to the
#   "date + time" format already used by B2:B4.
# - A new row 6 is appended, carrying the "date-only" format that B5 used
#   to have, together with new data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B5 adopts the same number format as B2:B4 (yyyy-mm-dd h:mm:ss style).
$ws.Range("B5").NumberFormat = $ws.Range("B4").NumberFormat

# New row 6: values for A6, B6, C6.
$ws.Range("A6").Value = 802.724
$ws.Range("B6").Value = 45728
$ws.Range("C6").Value = 810.465

# B6 keeps the number format that B5 used before this edit (yyyy-mm-dd).
$ws.Range("B6").NumberFormat = "YYYY-MM-DD"
